$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-19 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-20 Friday", 2) | Out-Null
$d.Content.Find.Execute("71-18=53", $true, $false, $false, $false, $false, $true, 1, $false, "71+1=72", 2) | Out-Null
$d.Content.Find.Execute("3+83=86", $true, $false, $false, $false, $false, $true, 1, $false, "85-13=72", 2) | Out-Null
$d.Content.Find.Execute("71-44=27", $true, $false, $false, $false, $false, $true, 1, $false, "80-58=22", 2) | Out-Null
$d.Content.Find.Execute("28-17=11", $true, $false, $false, $false, $false, $true, 1, $false, "88-28=60", 2) | Out-Null
$d.Content.Find.Execute("95-10=85", $true, $false, $false, $false, $false, $true, 1, $false, "31-17=14", 2) | Out-Null
$d.Content.Find.Execute("7+76=83", $true, $false, $false, $false, $false, $true, 1, $false, "80+4=84", 2) | Out-Null
$d.Content.Find.Execute("99-47=52", $true, $false, $false, $false, $false, $true, 1, $false, "23+68=91", 2) | Out-Null
$d.Content.Find.Execute("91-2=89", $true, $false, $false, $false, $false, $true, 1, $false, "49+8=57", 2) | Out-Null
$d.Content.Find.Execute("1+19=20", $true, $false, $false, $false, $false, $true, 1, $false, "94-5=89", 2) | Out-Null
$d.Content.Find.Execute("14+16=30", $true, $false, $false, $false, $false, $true, 1, $false, "77-58=19", 2) | Out-Null
$d.Content.Find.Execute("97-11=86", $true, $false, $false, $false, $false, $true, 1, $false, "52+30=82", 2) | Out-Null
$d.Content.Find.Execute("95-23=72", $true, $false, $false, $false, $false, $true, 1, $false, "96-44=52", 2) | Out-Null
$d.Content.Find.Execute("91+0=91", $true, $false, $false, $false, $false, $true, 1, $false, "70-47=23", 2) | Out-Null
$d.Content.Find.Execute("17+77=94", $true, $false, $false, $false, $false, $true, 1, $false, "78-18=60", 2) | Out-Null
$d.Content.Find.Execute("21-9=12", $true, $false, $false, $false, $false, $true, 1, $false, "42-16=26", 2) | Out-Null
$d.Content.Find.Execute("51-34=17", $true, $false, $false, $false, $false, $true, 1, $false, "13+13=26", 2) | Out-Null
$d.Content.Find.Execute("96-77=19", $true, $false, $false, $false, $false, $true, 1, $false, "54-7=47", 2) | Out-Null
$d.Content.Find.Execute("79+11=90", $true, $false, $false, $false, $false, $true, 1, $false, "70-38=32", 2) | Out-Null
$d.Content.Find.Execute("94-59=35", $true, $false, $false, $false, $false, $true, 1, $false, "77-69=8", 2) | Out-Null
$d.Content.Find.Execute("94-4=90", $true, $false, $false, $false, $false, $true, 1, $false, "38+44=82", 2) | Out-Null
$d.Content.Find.Execute("62-43=19", $true, $false, $false, $false, $false, $true, 1, $false, "73-40=33", 2) | Out-Null
$d.Content.Find.Execute("65-64=1", $true, $false, $false, $false, $false, $true, 1, $false, "37+53=90", 2) | Out-Null
$d.Content.Find.Execute("3+31=34", $true, $false, $false, $false, $false, $true, 1, $false, "6+88=94", 2) | Out-Null
$d.Content.Find.Execute("87-10=77", $true, $false, $false, $false, $false, $true, 1, $false, "46+51=97", 2) | Out-Null
$d.Content.Find.Execute("75-31=44", $true, $false, $false, $false, $false, $true, 1, $false, "84+7=91", 2) | Out-Null
$d.Content.Find.Execute("38+49=87", $true, $false, $false, $false, $false, $true, 1, $false, "20-17=3", 2) | Out-Null
$d.Content.Find.Execute("94-8=86", $true, $false, $false, $false, $false, $true, 1, $false, "86-55=31", 2) | Out-Null
$d.Content.Find.Execute("85-24=61", $true, $false, $false, $false, $false, $true, 1, $false, "67-41=26", 2) | Out-Null
$d.Content.Find.Execute("94-89=5", $true, $false, $false, $false, $false, $true, 1, $false, "17+71=88", 2) | Out-Null
$d.Content.Find.Execute("22+56=78", $true, $false, $false, $false, $false, $true, 1, $false, "12+61=73", 2) | Out-Null
$d.Content.Find.Execute("68-50=18", $true, $false, $false, $false, $false, $true, 1, $false, "4+23=27", 2) | Out-Null
$d.Content.Find.Execute("74-71=3", $true, $false, $false, $false, $false, $true, 1, $false, "0+84=84", 2) | Out-Null
$d.Content.Find.Execute("55-10=45", $true, $false, $false, $false, $false, $true, 1, $false, "57+13=70", 2) | Out-Null
$d.Content.Find.Execute("18+10=28", $true, $false, $false, $false, $false, $true, 1, $false, "75-2=73", 2) | Out-Null
$d.Content.Find.Execute("69-23=46", $true, $false, $false, $false, $false, $true, 1, $false, "70-41=29", 2) | Out-Null
$d.Content.Find.Execute("75-27=48", $true, $false, $false, $false, $false, $true, 1, $false, "13-5=8", 2) | Out-Null
$d.Content.Find.Execute("26+19=45", $true, $false, $false, $false, $false, $true, 1, $false, "66-52=14", 2) | Out-Null
$d.Content.Find.Execute("31-8=23", $true, $false, $false, $false, $false, $true, 1, $false, "27-22=5", 2) | Out-Null
$d.Content.Find.Execute("73-68=5", $true, $false, $false, $false, $false, $true, 1, $false, "19+0=19", 2) | Out-Null
$d.Content.Find.Execute("83-23=60", $true, $false, $false, $false, $false, $true, 1, $false, "69+8=77", 2) | Out-Null
$d.Content.Find.Execute("34+48=82", $true, $false, $false, $false, $false, $true, 1, $false, "89-53=36", 2) | Out-Null
$d.Content.Find.Execute("92-69=23", $true, $false, $false, $false, $false, $true, 1, $false, "43+56=99", 2) | Out-Null
$d.Content.Find.Execute("66+9=75", $true, $false, $false, $false, $false, $true, 1, $false, "59-13=46", 2) | Out-Null
$d.Content.Find.Execute("82-78=4", $true, $false, $false, $false, $false, $true, 1, $false, "59-7=52", 2) | Out-Null
$d.Content.Find.Execute("66-20=46", $true, $false, $false, $false, $false, $true, 1, $false, "17+59=76", 2) | Out-Null
$d.Content.Find.Execute("21+54=75", $true, $false, $false, $false, $false, $true, 1, $false, "65-55=10", 2) | Out-Null
$d.Content.Find.Execute("73-67=6", $true, $false, $false, $false, $false, $true, 1, $false, "32+6=38", 2) | Out-Null
$d.Content.Find.Execute("69-68=1", $true, $false, $false, $false, $false, $true, 1, $false, "14+1=15", 2) | Out-Null
$d.Content.Find.Execute("75-47=28", $true, $false, $false, $false, $false, $true, 1, $false, "26+25=51", 2) | Out-Null
$d.Content.Find.Execute("98-57=41", $true, $false, $false, $false, $false, $true, 1, $false, "27+1=28", 2) | Out-Null
$d.Content.Find.Execute("83-65=18", $true, $false, $false, $false, $false, $true, 1, $false, "30+64=94", 2) | Out-Null
$d.Content.Find.Execute("1+81=82", $true, $false, $false, $false, $false, $true, 1, $false, "24+22=46", 2) | Out-Null
$d.Content.Find.Execute("71-47=24", $true, $false, $false, $false, $false, $true, 1, $false, "26+69=95", 2) | Out-Null
$d.Content.Find.Execute("85-84=1", $true, $false, $false, $false, $false, $true, 1, $false, "9+41=50", 2) | Out-Null
$d.Content.Find.Execute("27-26=1", $true, $false, $false, $false, $false, $true, 1, $false, "29+56=85", 2) | Out-Null
$d.Content.Find.Execute("52-35=17", $true, $false, $false, $false, $false, $true, 1, $false, "6+38=44", 2) | Out-Null
$d.Content.Find.Execute("75-53=22", $true, $false, $false, $false, $false, $true, 1, $false, "97-4=93", 2) | Out-Null
$d.Content.Find.Execute("91-50=41", $true, $false, $false, $false, $false, $true, 1, $false, "25-7=18", 2) | Out-Null
$d.Content.Find.Execute("0+79=79", $true, $false, $false, $false, $false, $true, 1, $false, "13+80=93", 2) | Out-Null
$d.Content.Find.Execute("50+2=52", $true, $false, $false, $false, $false, $true, 1, $false, "53+19=72", 2) | Out-Null
$d.Content.Find.Execute("33+60=93", $true, $false, $false, $false, $false, $true, 1, $false, "83-66=17", 2) | Out-Null
$d.Content.Find.Execute("93-36=57", $true, $false, $false, $false, $false, $true, 1, $false, "25+16=41", 2) | Out-Null
$d.Content.Find.Execute("42+39=81", $true, $false, $false, $false, $false, $true, 1, $false, "26+66=92", 2) | Out-Null
$d.Content.Find.Execute("88+3=91", $true, $false, $false, $false, $false, $true, 1, $false, "51+23=74", 2) | Out-Null
$d.Content.Find.Execute("11+2=13", $true, $false, $false, $false, $false, $true, 1, $false, "52+32=84", 2) | Out-Null
$d.Content.Find.Execute("34-29=5", $true, $false, $false, $false, $false, $true, 1, $false, "95-41=54", 2) | Out-Null
$d.Content.Find.Execute("33+6=39", $true, $false, $false, $false, $false, $true, 1, $false, "74-27=47", 2) | Out-Null
$d.Content.Find.Execute("44+42=86", $true, $false, $false, $false, $false, $true, 1, $false, "98-87=11", 2) | Out-Null
$d.Content.Find.Execute("82-32=50", $true, $false, $false, $false, $false, $true, 1, $false, "98-30=68", 2) | Out-Null
$d.Content.Find.Execute("9+31=40", $true, $false, $false, $false, $false, $true, 1, $false, "15+40=55", 2) | Out-Null
$d.Content.Find.Execute("3+56=59", $true, $false, $false, $false, $false, $true, 1, $false, "4+48=52", 2) | Out-Null
$d.Content.Find.Execute("62-48=14", $true, $false, $false, $false, $false, $true, 1, $false, "56+8=64", 2) | Out-Null
$d.Content.Find.Execute("72-37=35", $true, $false, $false, $false, $false, $true, 1, $false, "11+69=80", 2) | Out-Null
$d.Content.Find.Execute("28+39=67", $true, $false, $false, $false, $false, $true, 1, $false, "79+7=86", 2) | Out-Null
$d.Content.Find.Execute("80-17=63", $true, $false, $false, $false, $false, $true, 1, $false, "16+77=93", 2) | Out-Null
$d.Content.Find.Execute("17+5=22", $true, $false, $false, $false, $false, $true, 1, $false, "40-7=33", 2) | Out-Null
$d.Content.Find.Execute("21+74=95", $true, $false, $false, $false, $false, $true, 1, $false, "5+1=6", 2) | Out-Null
$d.Content.Find.Execute("97-76=21", $true, $false, $false, $false, $false, $true, 1, $false, "89-32=57", 2) | Out-Null
$d.Content.Find.Execute("4+31=35", $true, $false, $false, $false, $false, $true, 1, $false, "26+65=91", 2) | Out-Null
$d.Content.Find.Execute("14+5=19", $true, $false, $false, $false, $false, $true, 1, $false, "87-54=33", 2) | Out-Null
$d.Content.Find.Execute("84-28=56", $true, $false, $false, $false, $false, $true, 1, $false, "50+49=99", 2) | Out-Null
$d.Content.Find.Execute("45+18=63", $true, $false, $false, $false, $false, $true, 1, $false, "53+3=56", 2) | Out-Null
$d.Content.Find.Execute("24+29=53", $true, $false, $false, $false, $false, $true, 1, $false, "29+28=57", 2) | Out-Null
$d.Content.Find.Execute("64-40=24", $true, $false, $false, $false, $false, $true, 1, $false, "11+79=90", 2) | Out-Null
$d.Content.Find.Execute("52-25=27", $true, $false, $false, $false, $false, $true, 1, $false, "25+73=98", 2) | Out-Null
$d.Content.Find.Execute("41-40=1", $true, $false, $false, $false, $false, $true, 1, $false, "31+59=90", 2) | Out-Null
$d.Content.Find.Execute("0+81=81", $true, $false, $false, $false, $false, $true, 1, $false, "57-24=33", 2) | Out-Null
$d.Content.Find.Execute("85-70=15", $true, $false, $false, $false, $false, $true, 1, $false, "46+49=95", 2) | Out-Null
$d.Content.Find.Execute("56+14=70", $true, $false, $false, $false, $false, $true, 1, $false, "11+22=33", 2) | Out-Null
$d.Content.Find.Execute("42+11=53", $true, $false, $false, $false, $false, $true, 1, $false, "90-87=3", 2) | Out-Null
$d.Content.Find.Execute("14+46=60", $true, $false, $false, $false, $false, $true, 1, $false, "59-24=35", 2) | Out-Null
$d.Content.Find.Execute("19+46=65", $true, $false, $false, $false, $false, $true, 1, $false, "23+62=85", 2) | Out-Null
$d.Content.Find.Execute("45+14=59", $true, $false, $false, $false, $false, $true, 1, $false, "92-16=76", 2) | Out-Null
$d.Content.Find.Execute("89-1=88", $true, $false, $false, $false, $false, $true, 1, $false, "31-28=3", 2) | Out-Null
$d.Content.Find.Execute("33+51=84", $true, $false, $false, $false, $false, $true, 1, $false, "2+25=27", 2) | Out-Null
$d.Content.Find.Execute("1+17=18", $true, $false, $false, $false, $false, $true, 1, $false, "4+64=68", 2) | Out-Null
$d.Content.Find.Execute("43-9=34", $true, $false, $false, $false, $false, $true, 1, $false, "49-37=12", 2) | Out-Null
$d.Content.Find.Execute("29-24=5", $true, $false, $false, $false, $false, $true, 1, $false, "5+45=50", 2) | Out-Null
$d.Content.Find.Execute("13+32=45", $true, $false, $false, $false, $false, $true, 1, $false, "21+25=46", 2) | Out-Null
$d.Content.Find.Execute("23-0=23", $true, $false, $false, $false, $false, $true, 1, $false, "83-27=56", 2) | Out-Null

Write-Host "Replacements complete"
